# Auto-generated: applies cryptos.xlsx data refresh (price/volume updates + row 25/26 swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.033.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.314.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.308.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.627"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.831.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.302.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.819.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.979"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.72%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "61.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "562.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.22%  "

$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -1.60%  "

$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0722"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.032.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0411"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("E47").Value = "  +2.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.61%  "

$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "
